$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) for columns A-E so that numeric-looking
# IDs and ISO dates are stored as literal text instead of being auto-
# converted to numbers/dates, matching the source data. The format is
# cleared again right after the value is written so the cells keep the
# default (unstyled) appearance.
$newRowsRange = $ws.Range("A108:E123")
$newRowsRange.NumberFormat = "@"

$ws.Range("A108").Value = "14357966"
$ws.Range("B108").Value = "2025-08-08"
$ws.Range("C108").Value = "Miomir Kecmanovic"
$ws.Range("D108").Value = "Ethan Quinn"
$ws.Range("E108").Value = "Gana Ethan Quinn"
$ws.Range("F108").Value = 2.1
$ws.Range("G107").Copy($ws.Range("G108"))
$ws.Range("H107").Copy($ws.Range("H108"))

$ws.Range("A109").Value = "14357997"
$ws.Range("B109").Value = "2025-08-08"
$ws.Range("C109").Value = "Hugo Dellien"
$ws.Range("D109").Value = "Reilly Opelka"
$ws.Range("E109").Value = "Gana Hugo Dellien"
$ws.Range("F109").Value = 6
$ws.Range("G107").Copy($ws.Range("G109"))
$ws.Range("H107").Copy($ws.Range("H109"))

$ws.Range("A110").Value = "14357999"
$ws.Range("B110").Value = "2025-08-08"
$ws.Range("C110").Value = "Jenson Brooksby"
$ws.Range("D110").Value = "Alexandre Muller"
$ws.Range("E110").Value = "Gana Alexandre Muller"
$ws.Range("F110").Value = 2.1
$ws.Range("G107").Copy($ws.Range("G110"))
$ws.Range("H107").Copy($ws.Range("H110"))

$ws.Range("A111").Value = "14357972"
$ws.Range("B111").Value = "2025-08-09"
$ws.Range("C111").Value = "Gael Monfils"
$ws.Range("D111").Value = "Nishesh Basavareddy"
$ws.Range("E111").Value = "Gana Nishesh Basavareddy"
$ws.Range("F111").Value = 2.3
$ws.Range("G107").Copy($ws.Range("G111"))
$ws.Range("H107").Copy($ws.Range("H111"))

$ws.Range("A112").Value = "14366985"
$ws.Range("B112").Value = "2025-08-08"
$ws.Range("C112").Value = "Caroline Dolehide"
$ws.Range("D112").Value = "Rebecca Sramkova"
$ws.Range("E112").Value = "Gana Rebecca Sramkova"
$ws.Range("F112").Value = 2.3
$ws.Range("G107").Copy($ws.Range("G112"))
$ws.Range("H107").Copy($ws.Range("H112"))

$ws.Range("A113").Value = "14369265"
$ws.Range("B113").Value = "2025-08-08"
$ws.Range("C113").Value = "Yuliia Starodubtseva"
$ws.Range("D113").Value = "Leolia Jeanjean"
$ws.Range("E113").Value = "Gana Leolia Jeanjean"
$ws.Range("F113").Value = 2.38
$ws.Range("G107").Copy($ws.Range("G113"))
$ws.Range("H107").Copy($ws.Range("H113"))

$ws.Range("A114").Value = "14366984"
$ws.Range("B114").Value = "2025-08-08"
$ws.Range("C114").Value = "Anna Blinkova"
$ws.Range("D114").Value = "Kimberly Birrell"
$ws.Range("E114").Value = "Gana Anna Blinkova"
$ws.Range("F114").Value = 1.8
$ws.Range("G107").Copy($ws.Range("G114"))
$ws.Range("H107").Copy($ws.Range("H114"))

$ws.Range("A115").Value = "14366987"
$ws.Range("B115").Value = "2025-08-08"
$ws.Range("C115").Value = "Anna Bondar"
$ws.Range("D115").Value = "Ajla Tomljanovic"
$ws.Range("E115").Value = "Gana Anna Bondar"
$ws.Range("F115").Value = 2.5
$ws.Range("G107").Copy($ws.Range("G115"))
$ws.Range("H107").Copy($ws.Range("H115"))

$ws.Range("A116").Value = "14369263"
$ws.Range("B116").Value = "2025-08-08"
$ws.Range("C116").Value = "Polina Kudermetova"
$ws.Range("D116").Value = "Ella Seidel"
$ws.Range("E116").Value = "Gana Polina Kudermetova"
$ws.Range("F116").Value = 2.2
$ws.Range("G107").Copy($ws.Range("G116"))
$ws.Range("H107").Copy($ws.Range("H116"))

$ws.Range("A117").Value = "14369264"
$ws.Range("B117").Value = "2025-08-08"
$ws.Range("C117").Value = "Kamilla Rakhimova"
$ws.Range("D117").Value = "Maria Sakkari"
$ws.Range("E117").Value = "Gana Kamilla Rakhimova"
$ws.Range("F117").Value = 3.2
$ws.Range("G107").Copy($ws.Range("G117"))
$ws.Range("H107").Copy($ws.Range("H117"))

$ws.Range("A118").Value = "14369266"
$ws.Range("B118").Value = "2025-08-08"
$ws.Range("C118").Value = "Anastasija Sevastova"
$ws.Range("D118").Value = "Emina Bektas"
$ws.Range("E118").Value = "Gana Emina Bektas"
$ws.Range("F118").Value = 2.75
$ws.Range("G107").Copy($ws.Range("G118"))
$ws.Range("H107").Copy($ws.Range("H118"))

$ws.Range("A119").Value = "14366986"
$ws.Range("B119").Value = "2025-08-08"
$ws.Range("C119").Value = "Suzan Lamens"
$ws.Range("D119").Value = "Veronika Kudermetova"
$ws.Range("E119").Value = "Gana Suzan Lamens"
$ws.Range("F119").Value = 3.75
$ws.Range("G107").Copy($ws.Range("G119"))
$ws.Range("H107").Copy($ws.Range("H119"))

$ws.Range("A120").Value = "14366983"
$ws.Range("B120").Value = "2025-08-08"
$ws.Range("C120").Value = "Alycia Parks"
$ws.Range("D120").Value = "Barbora Krejcikova"
$ws.Range("E120").Value = "Gana Alycia Parks"
$ws.Range("F120").Value = 2.75
$ws.Range("G107").Copy($ws.Range("G120"))
$ws.Range("H107").Copy($ws.Range("H120"))

$ws.Range("A121").Value = "14369261"
$ws.Range("B121").Value = "2025-08-09"
$ws.Range("C121").Value = "Catherine McNally"
$ws.Range("D121").Value = "Maddison Inglis"
$ws.Range("E121").Value = "Gana Maddison Inglis"
$ws.Range("F121").Value = 3.5
$ws.Range("G107").Copy($ws.Range("G121"))
$ws.Range("H107").Copy($ws.Range("H121"))

$ws.Range("A122").Value = "14311077"
$ws.Range("B122").Value = "2025-08-08"
$ws.Range("C122").Value = "Lukas Neumayer"
$ws.Range("D122").Value = "Nicolas Kicker"
$ws.Range("E122").Value = "Gana Nicolas Kicker"
$ws.Range("F122").Value = 2.63
$ws.Range("G107").Copy($ws.Range("G122"))
$ws.Range("H107").Copy($ws.Range("H122"))

$ws.Range("A123").Value = "14310258"
$ws.Range("B123").Value = "2025-08-08"
$ws.Range("C123").Value = "Yibing Wu"
$ws.Range("D123").Value = "Yu Hsiou Hsu"
$ws.Range("E123").Value = "Gana Yu Hsiou Hsu"
$ws.Range("F123").Value = 3.5
$ws.Range("G107").Copy($ws.Range("G123"))
$ws.Range("H107").Copy($ws.Range("H123"))

# Remove the temporary Text format applied above, restoring the default
# (General / unstyled) look while keeping the values already entered as text.
$newRowsRange.ClearFormats()

Write-Output "Added rows 108-123"